$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Matriz de Confusao (Treino)" values
$ws.Range("A6").Value = 16020
$ws.Range("B6").Value = 2
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 19093

# Update accuracy/precision/recall/F1 (train) values
$ws.Range("B35").Value = 0.9999145688575009
$ws.Range("B39").Value = 0.9998952605394082
$ws.Range("B43").Value = 0.9999476275269719
$ws.Range("B47").Value = 0.9999214433475608

# Update execution time text
$ws.Range("B49").Value = "0:01:28.510974"

# Add new rows for best parameters
$ws.Range("A51").Value = "Melhores Parâmetros"

$ws.Range("A52").Value = "solver"
$ws.Range("B52").Value = "liblinear"

$ws.Range("A53").Value = "penalty"
$ws.Range("B53").Value = "l1"

$ws.Range("A54").Value = "C"
$ws.Range("B54").Value = 10

$wb.Save()
